{"js": "// Merge the two runs of the first paragraph (\"Folder_\" + \"00_File_01\") into a\n// single run reading \"Folder_01_Subfolder01_File_01\", and relocate the\n// \"_GoBack\" bookmark (currently around the empty spot after \"Page 1/1\" in the\n// second paragraph) so it wraps the new text in the first paragraph instead.\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = body.paragraphs.items[0];\n\n// Replace the whole paragraph's text (both runs) with the new combined text.\nconst wholeRange = firstParagraph.getRange(\"Whole\");\nwholeRange.insertText(\"Folder_01_Subfolder01_File_01\", \"Replace\");\nawait context.sync();\n\n// Remove the existing \"_GoBack\" bookmark from its old location.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Re-insert \"_GoBack\" around the first paragraph's content (not the paragraph\n// mark), so bookmarkStart/bookmarkEnd both end up wrapping the new run.\nconst contentRange = firstParagraph.getRange(\"Content\");\ncontentRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Merge the two runs of the first paragraph (\"Folder_\" + \"00_File_01\") into a\n# single run reading \"Folder_01_Subfolder01_File_01\", and relocate the\n# \"_GoBack\" bookmark (currently around the empty spot after \"Page 1/1\" in the\n# second paragraph) so it wraps the new text in the first paragraph instead.\n\n$d = $word.ActiveDocument\n\n# Combine the two runs' text into a single run via Find & Replace (Word\n# merges the matched span into one run when the text is replaced).\n$null = $d.Content.Find.Execute(\"Folder_00_File_01\", $false, $false, $false, $false, $false, $true, 1, $false, \"Folder_01_Subfolder01_File_01\", 2)\n\n# Remove the existing \"_GoBack\" bookmark from its old location.\n$bookmark = $d.Bookmarks(\"_GoBack\")\n$bookmark.Delete()\n\n# Re-insert \"_GoBack\" around the first paragraph's content, excluding the\n# trailing paragraph mark, so bookmarkStart/bookmarkEnd both wrap the run.\n$firstParagraphRange = $d.Paragraphs.Item(1).Range\n$firstParagraphRange.MoveEnd(1, -1)\n$d.Bookmarks.Add(\"_GoBack\", $firstParagraphRange)\n"}
